$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 1 de Mayo de 2020 a las 13:22"

# Row 12
$ws.Cells.Item(12, 2).Value = 95646
$ws.Cells.Item(12, 3).Value = 1006
$ws.Cells.Item(12, 4).Value = 76318
$ws.Cells.Item(12, 5).Value = 13237
$ws.Cells.Item(12, 6).Value = 2899
$ws.Cells.Item(12, 7).Value = 63
$ws.Cells.Item(12, 8).Value = 6091

# Row 20
$ws.Cells.Item(20, 2).Value = 29705
$ws.Cells.Item(20, 3).Value = 119
$ws.Cells.Item(20, 4).Value = 23400
$ws.Cells.Item(20, 5).Value = 4568
$ws.Cells.Item(20, 6).Value = 167

# Row 27
$ws.Cells.Item(27, 1).Value = "Pakistan"
$ws.Cells.Item(27, 2).Value = 17439
$ws.Cells.Item(27, 3).Value = 966
$ws.Cells.Item(27, 4).Value = 4315
$ws.Cells.Item(27, 5).Value = 12733
$ws.Cells.Item(27, 6).Value = 111
$ws.Cells.Item(27, 7).Value = 30
$ws.Cells.Item(27, 8).Value = 391

# Row 28
$ws.Cells.Item(28, 1).Value = "Singapur"
$ws.Cells.Item(28, 2).Value = 17101
$ws.Cells.Item(28, 3).Value = 932
$ws.Cells.Item(28, 4).Value = 1244
$ws.Cells.Item(28, 5).Value = 15842
$ws.Cells.Item(28, 6).Value = 21
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 15

# Row 33
$ws.Cells.Item(33, 1).Value = "Catar"
$ws.Cells.Item(33, 2).Value = 14096
$ws.Cells.Item(33, 3).Value = 687
$ws.Cells.Item(33, 4).Value = 1436
$ws.Cells.Item(33, 5).Value = 12648
$ws.Cells.Item(33, 6).Value = 72
$ws.Cells.Item(33, 7).Value = 2
$ws.Cells.Item(33, 8).Value = 12

# Row 34
$ws.Cells.Item(34, 1).Value = "Japon"
$ws.Cells.Item(34, 2).Value = 14088
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 2460
$ws.Cells.Item(34, 5).Value = 11198
$ws.Cells.Item(34, 6).Value = 308
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 430

# Row 35
$ws.Cells.Item(35, 1).Value = "Polonia"
$ws.Cells.Item(35, 2).Value = 13105
$ws.Cells.Item(35, 3).Value = 228
$ws.Cells.Item(35, 4).Value = 3491
$ws.Cells.Item(35, 5).Value = 8963
$ws.Cells.Item(35, 6).Value = 160
$ws.Cells.Item(35, 7).Value = 7
$ws.Cells.Item(35, 8).Value = 651

# Row 36
$ws.Cells.Item(36, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(36, 2).Value = 13038
$ws.Cells.Item(36, 3).Value = 557
$ws.Cells.Item(36, 4).Value = 2543
$ws.Cells.Item(36, 5).Value = 10384
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 6
$ws.Cells.Item(36, 8).Value = 111

# Row 48
$ws.Cells.Item(48, 2).Value = 6767
$ws.Cells.Item(48, 3).Value = 13
$ws.Cells.Item(48, 4).Value = 5745
$ws.Cells.Item(48, 5).Value = 929
$ws.Cells.Item(48, 6).Value = 28
$ws.Cells.Item(48, 7).Value = 1

# Row 54
$ws.Cells.Item(54, 5).Value = 1833
$ws.Cells.Item(54, 6).Value = 49
$ws.Cells.Item(54, 7).Value = 7
$ws.Cells.Item(54, 8).Value = 218

# Row 59
$ws.Cells.Item(59, 5).Value = 2503
$ws.Cells.Item(59, 6).Value = 237
$ws.Cells.Item(59, 7).Value = 6
$ws.Cells.Item(59, 8).Value = 122

# Row 71
$ws.Cells.Item(71, 1).Value = "Uzbekistan"
$ws.Cells.Item(71, 2).Value = 2075
$ws.Cells.Item(71, 3).Value = 36
$ws.Cells.Item(71, 4).Value = 1182
$ws.Cells.Item(71, 5).Value = 884
$ws.Cells.Item(71, 6).Value = 8
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 9

# Row 72
$ws.Cells.Item(72, 1).Value = "Ghana"
$ws.Cells.Item(72, 2).Value = 2074
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 212
$ws.Cells.Item(72, 5).Value = 1845
$ws.Cells.Item(72, 6).Value = 4
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 17

# Row 136
$ws.Cells.Item(136, 2).Value = 179
$ws.Cells.Item(136, 3).Value = 1
$ws.Cells.Item(136, 4).Value = 83
$ws.Cells.Item(136, 5).Value = 82
$ws.Cells.Item(136, 6).Value = 5

# Row 137
$ws.Cells.Item(137, 1).Value = "Guadalupe"
$ws.Cells.Item(137, 2).Value = 152
$ws.Cells.Item(137, 3).Value = 1
$ws.Cells.Item(137, 4).Value = 95
$ws.Cells.Item(137, 5).Value = 45
$ws.Cells.Item(137, 6).Value = 6
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 12

# Row 138
$ws.Cells.Item(138, 1).Value = "Birmania"
$ws.Cells.Item(138, 3).Value = 1
$ws.Cells.Item(138, 4).Value = 28
$ws.Cells.Item(138, 5).Value = 117
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 6

# Row 169
$ws.Cells.Item(169, 4).Value = 51
$ws.Cells.Item(169, 5).Value = 7
$ws.Cells.Item(169, 6).Value = 1
